# Delete row 666 ("ようこそ" post) entirely.
# This shifts all subsequent rows (667-673) up by one, matching the
# target diff where row 666 is removed and rows 667-673 become 666-672.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("666").Delete()
